$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 17
$ws.Range("H2").Value = 17
$ws.Range("E3").Value = 14
$ws.Range("F3").Value = 8
$ws.Range("H3").Value = 8
$ws.Range("E9").Value = 12
$ws.Range("E15").Value = 116
$ws.Range("F15").Value = 53
$ws.Range("H15").Value = 53
$ws.Range("E16").Value = 7
$ws.Range("F16").Value = 3
$ws.Range("H16").Value = 3
$ws.Range("E17").Value = 73
$ws.Range("F17").Value = 27
$ws.Range("H17").Value = 27
$ws.Range("E18").Value = 68
$ws.Range("F18").Value = 24
$ws.Range("H18").Value = 24
$ws.Range("E19").Value = 31
$ws.Range("F19").Value = 18
$ws.Range("H19").Value = 18
$ws.Range("F24").Value = 8
$ws.Range("H24").Value = 8
$ws.Range("E29").Value = 10
$ws.Range("F29").Value = 4
$ws.Range("H29").Value = 4
$ws.Range("E34").Value = 11
$ws.Range("E36").Value = 59
$ws.Range("F36").Value = 20
$ws.Range("H36").Value = 20
$ws.Range("E37").Value = 29
$ws.Range("E38").Value = 40
$ws.Range("E40").Value = 10
$ws.Range("E42").Value = 24
$ws.Range("E43").Value = 14
$ws.Range("E47").Value = 41
$ws.Range("E48").Value = 17
$ws.Range("E49").Value = 45
$ws.Range("F49").Value = 23
$ws.Range("H49").Value = 23
$ws.Range("E55").Value = 3
$ws.Range("E58").Value = 3
$ws.Range("E60").Value = 11
$ws.Range("E63").Value = 15
$ws.Range("E65").Value = 20
$ws.Range("F65").Value = 7
$ws.Range("H65").Value = 7
$ws.Range("E67").Value = 25
$ws.Range("E70").Value = 21
$ws.Range("E71").Value = 19
$ws.Range("E72").Value = 24
$ws.Range("F72").Value = 13
$ws.Range("H72").Value = 13
$ws.Range("E76").Value = 29
$ws.Range("F77").Value = 9
$ws.Range("H77").Value = 9
$ws.Range("E78").Value = 17
$ws.Range("F78").Value = 6
$ws.Range("H78").Value = 6
$ws.Range("E79").Value = 16
$ws.Range("E80").Value = 15
$ws.Range("F80").Value = 6
$ws.Range("H80").Value = 6
$ws.Range("E85").Value = 3
$ws.Range("E87").Value = 5
$ws.Range("F88").Value = 7
$ws.Range("H88").Value = 7
$ws.Range("E89").Value = 19